$wb = $excel.ActiveWorkbook

# --- "Hoja1": refresh the daily conversion summary text (A1) ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$resumen = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.87 = 41474.49 pesos`n✅ 41474.49 pesos = 9.82 = 962.07 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $resumen

# --- "tasas": refresh the automatically-updated rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 101.35
$ws2.Range("O10").Value = 4203.44
$ws2.Range("N12").Value = 4223
$ws2.Range("O12").Value = 97.95999999999999
